# Update FSAR plot data with infilled spawning escapement time series
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 4342.234444444444
$ws.Range("D2").Value2 = 0.627988425
$ws.Range("B3").Value2 = 9873.310495235024
$ws.Range("D3").Value2 = 0.768276334
$ws.Range("B4").Value2 = 7757.792063492064
$ws.Range("D4").Value2 = 0.735717679
$ws.Range("B5").Value2 = 10201.78173185951
$ws.Range("D5").Value2 = 0.65244786
$ws.Range("B6").Value2 = 5661.329557713053
$ws.Range("D6").Value2 = 0.757651465
$ws.Range("B7").Value2 = 7454.889411764706
$ws.Range("D7").Value2 = 0.75568565
$ws.Range("B8").Value2 = 6242.570323464136
$ws.Range("D8").Value2 = 0.679425549
$ws.Range("B9").Value2 = 6343.314175234535
$ws.Range("D9").Value2 = 0.643201701
$ws.Range("B10").Value2 = 4381.278666514841
$ws.Range("D10").Value2 = 0.261422611
$ws.Range("B11").Value2 = 8304.15571645675
$ws.Range("D11").Value2 = 0.347515901
$ws.Range("B12").Value2 = 9538.829496623695
$ws.Range("D12").Value2 = 0.28997512
$ws.Range("B13").Value2 = 8246.923787599755
$ws.Range("D13").Value2 = 0.347183989
$ws.Range("B14").Value2 = 8648.799346286722
$ws.Range("D14").Value2 = 0.491070568
$ws.Range("B15").Value2 = 10161.95270440327
$ws.Range("D15").Value2 = 0.639342583
$ws.Range("B16").Value2 = 11291.19643953346
$ws.Range("D16").Value2 = 0.594688916
$ws.Range("B17").Value2 = 11938.56384284837
$ws.Range("D17").Value2 = 0.640580262
$ws.Range("B18").Value2 = 6353.0
$ws.Range("D18").Value2 = 0.584398513
$ws.Range("B19").Value2 = 10609.0
$ws.Range("D19").Value2 = 0.217918697
$ws.Range("B20").Value2 = 14176.0
$ws.Range("D20").Value2 = 0.247923591
$ws.Range("B21").Value2 = 29840.0
$ws.Range("D21").Value2 = 0.376006231
$ws.Range("B22").Value2 = 20847.0
$ws.Range("D22").Value2 = 0.563086753
$ws.Range("B23").Value2 = 8110.0
$ws.Range("D23").Value2 = 0.386694006
$ws.Range("B24").Value2 = 8565.27380952381
$ws.Range("D24").Value2 = 0.070410852
$ws.Range("B25").Value2 = 22123.0
$ws.Range("D25").Value2 = 0.213216634
$ws.Range("B26").Value2 = 26789.0
$ws.Range("D26").Value2 = 0.214938854
$ws.Range("B27").Value2 = 23272.0
$ws.Range("D27").Value2 = 0.345757444
$ws.Range("B28").Value2 = 8608.999539594843
$ws.Range("D28").Value2 = 0.404666051
$ws.Range("B29").Value2 = 13582.08695652174
$ws.Range("D29").Value2 = 0.3429582
$ws.Range("B30").Value2 = 7977.826086956522
$ws.Range("D30").Value2 = 0.445921351
$ws.Range("B31").Value2 = 6683.0
$ws.Range("D31").Value2 = 0.351464235
$ws.Range("B32").Value2 = 13856.95652173913
$ws.Range("D32").Value2 = 0.483432028
$ws.Range("B33").Value2 = 10932.0
$ws.Range("D33").Value2 = 0.300333805
$ws.Range("B34").Value2 = 11295.0
$ws.Range("D34").Value2 = 0.345180591
$ws.Range("B35").Value2 = 8652.952437574317
$ws.Range("D35").Value2 = 0.430712842
$ws.Range("B36").Value2 = 17803.0
$ws.Range("D36").Value2 = 0.328617006
$ws.Range("B37").Value2 = 13288.0
$ws.Range("D37").Value2 = 0.440206153
$ws.Range("B38").Value2 = 23619.47086801427
$ws.Range("D38").Value2 = 0.255179599
$ws.Range("B39").Value2 = 23162.0
$ws.Range("D39").Value2 = 0.383416369
$ws.Range("B40").Value2 = 18250.0
$ws.Range("D40").Value2 = 0.439626851
$ws.Range("B41").Value2 = 13737.0
$ws.Range("D41").Value2 = 0.302697746
$ws.Range("B42").Value2 = 16721.62137578169
$ws.Range("D42").Value2 = 0.366007267
$ws.Range("B43").Value2 = 20125.0
$ws.Range("D43").Value2 = 0.280841141
$ws.Range("B44").Value2 = 21689.0
$ws.Range("D44").Value2 = 0.319950959
$ws.Range("B45").Value2 = 18243.0
$ws.Range("C45").Value2 = 160604.8106129762
$ws.Range("D45").Value2 = 0.406107345

# Append new row 46 (year 2023)
$ws.Range("A46").Value2 = 2023
$ws.Range("B46").Value2 = 33612.24729960204
$ws.Range("C46").Value2 = 83986.66795893584
$ws.Range("D46").Value2 = 0.234145889

# Restore the scrolled/selected view state recorded in the saved workbook
$ws.Range("F55").Select()
